$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: append a new "copied" table (with format) at rows 16-22 ---

# Header row 16 (A..E) - values first, then copy the format of the existing
# header block (A5:B5) over it so it picks up the same style (yellow fill +
# border + centered) without creating a new style entry.
$ws1.Range("A16").Value2 = "A"
$ws1.Range("B16").Value2 = "B"
$ws1.Range("C16").Value2 = "C"
$ws1.Range("D16").Value2 = "D"
$ws1.Range("E16").Value2 = "E"

$hdrSrc1 = $ws1.Range("A5:B5")
$hdrSrc1.Copy()
$ws1.Range("A16:E16").PasteSpecial(-4122)

# Data rows 17-22. Values are written as text formulas ("=""n""") so that,
# once flattened to static values below, they land as text (matching the
# destination format used elsewhere in this sheet) rather than numbers.
$ws1.Range("A17").Formula = "=""1"""
$ws1.Range("B17").Formula = "=""2"""
$ws1.Range("C17").Formula = "=""3"""
$ws1.Range("D17").Formula = "=""4"""
$ws1.Range("E17").Formula = "=""5"""

$ws1.Range("A18").Formula = "=""2"""
$ws1.Range("B18").Formula = "=""3"""
$ws1.Range("C18").Formula = "=""4"""
$ws1.Range("D18").Formula = "=""5"""
$ws1.Range("E18").Formula = "=""6"""

$ws1.Range("A19").Formula = "=""3"""
$ws1.Range("B19").Formula = "=""4"""
$ws1.Range("C19").Formula = "=""5"""
$ws1.Range("D19").Formula = "=""6"""
$ws1.Range("E19").Formula = "=""7"""

$ws1.Range("A20").Formula = "=""4"""
$ws1.Range("B20").Formula = "=""5"""
$ws1.Range("C20").Formula = "=""6"""
$ws1.Range("D20").Formula = "=""7"""
$ws1.Range("E20").Formula = "=""8"""

$ws1.Range("A21").Formula = "=""5"""
$ws1.Range("B21").Formula = "=""6"""
$ws1.Range("C21").Formula = "=""7"""
$ws1.Range("D21").Formula = "=""8"""
$ws1.Range("E21").Formula = "=""9"""

$ws1.Range("A22").Formula = "=""6"""
$ws1.Range("B22").Formula = "=""7"""
$ws1.Range("C22").Formula = "=""8"""
$ws1.Range("D22").Formula = "=""9"""
$ws1.Range("E22").Formula = "=""10"""

# Apply the border-only look used for data rows across the sheet (reuses
# the existing thin/black border definition).
$dataRng1 = $ws1.Range("A17:E22")
$dataRng1.Borders.Color = 0
$dataRng1.Borders.Weight = 2

# Flatten the text formulas down to static shared-string values in place.
$dataRng1.Copy()
$dataRng1.PasteSpecial(-4163)

# --- Add Sheet3 as a new destination sheet and copy the range there ---

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet3"

# Header row, starting at destination cell C5.
$ws3.Range("C5").Value2 = "A"
$ws3.Range("D5").Value2 = "B"
$ws3.Range("E5").Value2 = "C"
$ws3.Range("F5").Value2 = "D"
$ws3.Range("G5").Value2 = "E"

$hdrSrc3 = $ws1.Range("A5:B5")
$hdrSrc3.Copy()
$ws3.Range("C5:G5").PasteSpecial(-4122)

# Data rows, as plain numbers this time.
$ws3.Range("C6").Value2 = 1
$ws3.Range("D6").Value2 = 2
$ws3.Range("E6").Value2 = 3
$ws3.Range("F6").Value2 = 4
$ws3.Range("G6").Value2 = 5

$ws3.Range("C7").Value2 = 2
$ws3.Range("D7").Value2 = 3
$ws3.Range("E7").Value2 = 4
$ws3.Range("F7").Value2 = 5
$ws3.Range("G7").Value2 = 6

$ws3.Range("C8").Value2 = 3
$ws3.Range("D8").Value2 = 4
$ws3.Range("E8").Value2 = 5
$ws3.Range("F8").Value2 = 6
$ws3.Range("G8").Value2 = 7

$ws3.Range("C9").Value2 = 4
$ws3.Range("D9").Value2 = 5
$ws3.Range("E9").Value2 = 6
$ws3.Range("F9").Value2 = 7
$ws3.Range("G9").Value2 = 8

$ws3.Range("C10").Value2 = 5
$ws3.Range("D10").Value2 = 6
$ws3.Range("E10").Value2 = 7
$ws3.Range("F10").Value2 = 8
$ws3.Range("G10").Value2 = 9

$ws3.Range("C11").Value2 = 6
$ws3.Range("D11").Value2 = 7
$ws3.Range("E11").Value2 = 8
$ws3.Range("F11").Value2 = 9
$ws3.Range("G11").Value2 = 10

$dataRng3 = $ws3.Range("C6:G11")
$dataRng3.Borders.Color = 0
$dataRng3.Borders.Weight = 2

# Leave the destination cell selected on the new sheet, which becomes the
# active sheet/tab.
$ws3.Range("J11").Select()
